$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "KEY ACHIEVEMENTS AND IMPACT" section so we only touch
# the bullet paragraphs that live there (several similar / duplicate
# sentences exist elsewhere in the resume, e.g. under PROFESSIONAL
# EXPERIENCE and KEY PROJECTS).
# ------------------------------------------------------------------
$sectionHeading = "KEY ACHIEVEMENTS AND IMPACT"
$count = $d.Paragraphs.Count
$sectionStart = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq $sectionHeading) {
        $sectionStart = $i
        break
    }
}

if ($sectionStart -eq -1) {
    throw "Could not find section heading '$sectionHeading'"
}

# Find the end of the section (next Heading2-styled paragraph, or end
# of document) so we know the bounds we are allowed to edit within.
$sectionEnd = $count
for ($i = $sectionStart + 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Heading 2") {
        $sectionEnd = $i - 1
        break
    }
}

# ------------------------------------------------------------------
# Helper: find, within [sectionStart, sectionEnd], the paragraph whose
# text contains the given needle, and return its index.
# ------------------------------------------------------------------
function Find-ParagraphIndex($needle) {
    for ($i = $sectionStart; $i -le $sectionEnd; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Replace three bullets in place (text-only change, keeps the
#    paragraph / run formatting and bullet prefix intact).
# ------------------------------------------------------------------
$idx = Find-ParagraphIndex("Built real-time FEC analysis systems using Python, Pandas and PySpark")
if ($idx -eq -1) { throw "Could not find FEC analysis bullet" }
$d.Paragraphs.Item($idx).Range.Find.Execute(
    "Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%", 2)

$idx = Find-ParagraphIndex("Built cloud-based data warehouse solutions on AWS")
if ($idx -eq -1) { throw "Could not find cloud-based data warehouse bullet" }
$d.Paragraphs.Item($idx).Range.Find.Execute(
    "Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access", 2)

$idx = Find-ParagraphIndex("Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS")
if ($idx -eq -1) { throw "Could not find ETL pipelines bullet" }
$d.Paragraphs.Item($idx).Range.Find.Execute(
    "Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions", 2)

# ------------------------------------------------------------------
# 2) Replace the final "race coding errors" bullet BEFORE deleting
#    the two bullets that sit between it and the earlier ones, so its
#    paragraph index doesn't shift under us.
# ------------------------------------------------------------------
$idx = Find-ParagraphIndex("Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%")
if ($idx -eq -1) { throw "Could not find race coding errors bullet" }
$d.Paragraphs.Item($idx).Range.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "178% accuracy improvement in racial classification algorithms", 2)

# ------------------------------------------------------------------
# 3) Delete the two bullets that are dropped entirely. Delete the
#    higher-indexed one first so the lower index stays valid.
# ------------------------------------------------------------------
$idx2 = Find-ParagraphIndex("Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations")
if ($idx2 -eq -1) { throw "Could not find redistricting platform bullet" }

$idx1 = Find-ParagraphIndex("Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis")
if ($idx1 -eq -1) { throw "Could not find trigonometric algorithm bullet" }

if ($idx2 -gt $idx1) {
    $d.Paragraphs.Item($idx2).Range.Delete()
    $d.Paragraphs.Item($idx1).Range.Delete()
} else {
    $d.Paragraphs.Item($idx1).Range.Delete()
    $d.Paragraphs.Item($idx2).Range.Delete()
}

Write-Host "Done editing KEY ACHIEVEMENTS AND IMPACT section."
